$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSFT")

$ws.Range("B4").Value  = 1916000000.0
$ws.Range("G4").Value  = 1823000000.0

$ws.Range("B7").Value  = 65618000000.0
$ws.Range("G7").Value  = 40522000000.0

$ws.Range("B8").Value  = 5429000000.0

$ws.Range("B10").Value = 14393000000.0
$ws.Range("G10").Value = 14630000000.0

$ws.Range("B13").Value = 14245000000.0
$ws.Range("G13").Value = 8811000000.0

$ws.Range("B15").Value = 10680000000.0
$ws.Range("G15").Value = 6247000000.0

$ws.Range("B18").Value = 45936000000.0
$ws.Range("G18").Value = 9131000000.0

$ws.Range("B20").Value = 70580000000.0
$ws.Range("G20").Value = 63361000000.0

$ws.Range("B22").Value = 85000000.0
$ws.Range("G22").Value = 222000000.0

$ws.Range("B23").Value = 28797000000.0
$ws.Range("G23").Value = 38412000000.0

$ws.Range("B33").Value = -44147000000.0

$ws.Range("B34").Value = 81260000000.0
